# Generate Report for Handback
# ------------------------------------------------------------------
# This mirrors the "localization-status.xlsx" handback report refresh:
#   * Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#     (Overview sheet E/F columns, and the Status column on each language sheet)
#   * Each language sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated now that a handback
#     round-trip has happened:
#       - Latest Target File: hyperlinked .md source file (same link target
#         as the "Source File Name" column)
#       - Latest Handback File: the generated target-language .xlf file name
#       - Latest Handback DateTime: the handback timestamp
#   * Columns that now hold longer text are widened to fit.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8af89d09c7342875bc033fa56f29189d92947dbf/e2e/"

# ------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells (E = zh-cn, F = de-de)
# ------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Columns.Item(5).ColumnWidth = 29.1
$ov.Columns.Item(6).ColumnWidth = 29.1

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (C)
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

# Row 2 -> 25cb016f-0214-4cbd-a58e-ba76f3905fe1
$zh.Hyperlinks.Add($zh.Range("I2"), $repoBase + "25cb016f-0214-4cbd-a58e-ba76f3905fe1.md", "", "", "25cb016f-0214-4cbd-a58e-ba76f3905fe1.md")
$zh.Range("J2").Value = "25cb016f-0214-4cbd-a58e-ba76f3905fe1.8486084da1a92d7c272a56f00b6ca422a2fe2065.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-03 08:51:40"

# Row 3 -> 30e7c031-6d43-4cbf-85f7-b48728048a7a
$zh.Hyperlinks.Add($zh.Range("I3"), $repoBase + "30e7c031-6d43-4cbf-85f7-b48728048a7a.md", "", "", "30e7c031-6d43-4cbf-85f7-b48728048a7a.md")
$zh.Range("J3").Value = "30e7c031-6d43-4cbf-85f7-b48728048a7a.453008b0ba8d0fbc233f6f842f73620b0c21786f.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-03 08:51:40"

# Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
$zh.Columns.Item(3).ColumnWidth = 29.1
$zh.Columns.Item(9).ColumnWidth = 39.1
$zh.Columns.Item(10).ColumnWidth = 39.1

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Status column (C)
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# Row 2 -> 25cb016f-0214-4cbd-a58e-ba76f3905fe1
$de.Hyperlinks.Add($de.Range("I2"), $repoBase + "25cb016f-0214-4cbd-a58e-ba76f3905fe1.md", "", "", "25cb016f-0214-4cbd-a58e-ba76f3905fe1.md")
$de.Range("J2").Value = "25cb016f-0214-4cbd-a58e-ba76f3905fe1.8486084da1a92d7c272a56f00b6ca422a2fe2065.de-de.xlf"
$de.Range("K2").Value = "2016-09-03 08:51:47"

# Row 3 -> 30e7c031-6d43-4cbf-85f7-b48728048a7a
$de.Hyperlinks.Add($de.Range("I3"), $repoBase + "30e7c031-6d43-4cbf-85f7-b48728048a7a.md", "", "", "30e7c031-6d43-4cbf-85f7-b48728048a7a.md")
$de.Range("J3").Value = "30e7c031-6d43-4cbf-85f7-b48728048a7a.453008b0ba8d0fbc233f6f842f73620b0c21786f.de-de.xlf"
$de.Range("K3").Value = "2016-09-03 08:51:47"

# Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
$de.Columns.Item(3).ColumnWidth = 29.1
$de.Columns.Item(9).ColumnWidth = 39.1
$de.Columns.Item(10).ColumnWidth = 39.1

Write-Output "Handback report generated."
